$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.720.68'
$ws.Range("E2").Value = '  +2.37%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.660.32'
$ws.Range("E3").Value = '  +2.24%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.75'
$ws.Range("E5").Value = '  +1.38%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.59'
$ws.Range("E6").Value = '  +1.82%  '

# Row 7
$ws.Range("E7").Value = '  -0.04%  '

# Row 8
$ws.Range("E8").Value = '  +0.74%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.60'
$ws.Range("E9").Value = '  +0.51%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.111'
$ws.Range("E10").Value = '  +4.68%  '

# Row 11
$ws.Range("E11").Value = '  +3.51%  '

# Row 12
$ws.Range("E12").Value = '  +0.84%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.129.65'
$ws.Range("E13").Value = '  +2.38%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.17'
$ws.Range("E14").Value = '  +7.49%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '61.591.77'
$ws.Range("E15").Value = '  +2.17%  '

# Row 16
$ws.Range("E16").Value = '  +4.68%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.672.13'
$ws.Range("E17").Value = '  +2.61%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.69'
$ws.Range("E18").Value = '  +2.88%  '

# Row 19
$ws.Range("E19").Value = '  +3.98%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '356.83'
$ws.Range("E20").Value = '  +3.06%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.92'
$ws.Range("E21").Value = '  +0.52%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.08%  '

# Row 23
$ws.Range("E23").Value = '  -0.78%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '64.85'
$ws.Range("E24").Value = '  +2.88%  '

# Row 25
$ws.Range("E25").Value = '  +3.45%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.996'
$ws.Range("E26").Value = '  -0.10%  '

# Row 27
$ws.Range("E27").Value = '  +5.57%  '

# Row 28
$ws.Range("E28").Value = '  +8.13%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0831'
$ws.Range("E29").Value = '  +4.24%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.97'
$ws.Range("E30").Value = '  +9.38%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '170.05'
$ws.Range("E31").Value = '  +2.27%  '

# Row 32
$ws.Range("E32").Value = '  -0.03%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.21'
$ws.Range("E33").Value = '  +4.11%  '

# Row 34
$ws.Range("E34").Value = '  +15.67%  '

# Row 35
$ws.Range("E35").Value = '  +9.50%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.39'
$ws.Range("E36").Value = '  +6.95%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +19.68%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.75'
$ws.Range("E38").Value = '  +7.75%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '343.99'
$ws.Range("E39").Value = '  +9.78%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.17'
$ws.Range("E40").Value = '  +7.23%  '

# Row 41
$ws.Range("E41").Value = '  +1.57%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.47'
$ws.Range("E42").Value = '  +9.61%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.97'
$ws.Range("E43").Value = '  +5.72%  '

# Row 44
$ws.Range("E44").Value = '  +5.55%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '21.29'
$ws.Range("E45").Value = '  +6.58%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0256'
$ws.Range("E46").Value = '  +6.08%  '

# Row 47
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.631'
$ws.Range("E47").Value = '  +4.26%  '

# Row 48
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.89'
$ws.Range("E48").Value = '  +0.61%  '

# Row 49
$ws.Range("E49").Value = '  +0.98%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.996'
$ws.Range("E50").Value = '  -0.26%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.117.05'
$ws.Range("E51").Value = '  +5.18%  '
